$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 <-> Row 3 swap (Perfect Binding <-> Fold)
$ws.Range("B2").Value = "Fold"
$ws.Range("D2").Value = "Duplo-Collator"
$ws.Range("G2").Value = "1,020"
$ws.Range("L2").Value = ""
$ws.Range("O2").Value = "Duplo-Collator"
$ws.Range("P2").Value = "Duplo-Collator"

$ws.Range("B3").Value = "Perfect Binding"
$ws.Range("D3").Value = "Bourg Perfect Binder"
$ws.Range("G3").Value = "1,010"
$ws.Range("L3").Value = "ISM Chicago"
$ws.Range("O3").Value = "Bourg Perfect Binder"
$ws.Range("P3").Value = "Bourg Perfect Binder"

# Rows 5-9 cyclic rotation (Cover 4p group), plus "F 4x0" -> "F 5x0" text fix
$ws.Range("B5").Value = "Sheet-fed Press F 5x0"
$ws.Range("D5").Value = "S1 HEI 640C"
$ws.Range("G5").Value = "1,117"
$ws.Range("O5").Value = "S1 HEI 640C"
$ws.Range("P5").Value = "S1 HEI 640C"

$ws.Range("B6").Value = "Sheet-fed Press B 0x4"
$ws.Range("D6").Value = "S1 HEI 640C"
$ws.Range("G6").Value = "672"
$ws.Range("O6").Value = "S1 HEI 640C"
$ws.Range("P6").Value = "S1 HEI 640C"

$ws.Range("B7").Value = "Cut"
$ws.Range("D7").Value = "Cutter 45"""
$ws.Range("G7").Value = "640"
$ws.Range("O7").Value = "Cutter 45"""
$ws.Range("P7").Value = "Cutter 45"""

$ws.Range("B8").Value = "Lamination"
$ws.Range("D8").Value = "Laminator"
$ws.Range("G8").Value = "1,082"
$ws.Range("O8").Value = "Laminator"
$ws.Range("P8").Value = "Laminator"

$ws.Range("B9").Value = "Cut"
$ws.Range("D9").Value = "Cutter 45"""
$ws.Range("G9").Value = "510"
$ws.Range("O9").Value = "Cutter 45"""
$ws.Range("P9").Value = "Cutter 45"""

# Rows 13 <-> 14 swap (OKTP <-> Plate burn), plus PlannedQty 8 -> 9
$ws.Range("B13").Value = "Plate burn"
$ws.Range("D13").Value = "Plate Making"
$ws.Range("G13").Value = "9"
$ws.Range("O13").Value = "Plate Making"
$ws.Range("P13").Value = "Plate Making"

$ws.Range("B14").Value = "OKTP"
$ws.Range("D14").Value = "Ok to Plate"
$ws.Range("G14").Value = "9"
$ws.Range("O14").Value = "OKTP"
$ws.Range("P14").Value = "Ok to Plate`nIntegration Validation"
